# Update loading_percent values for rows 2-25 (data rows 0-23), columns B,C,D,F,G,I,J
# per commit "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (index 0)
$ws.Range("B2").Value = 18.91727032928506
$ws.Range("C2").Value = 16.09710850822301
$ws.Range("D2").Value = 7.879781859840354
$ws.Range("F2").Value = 45.18703441539783
$ws.Range("G2").Value = 3.694307544394126
$ws.Range("I2").Value = 28.2426519515638
$ws.Range("J2").Value = 11.47564482662765

# Row 3 (index 1)
$ws.Range("B3").Value = 18.42430170994087
$ws.Range("C3").Value = 15.58282759897601
$ws.Range("D3").Value = 7.872716041014177
$ws.Range("F3").Value = 44.89350579763195
$ws.Range("G3").Value = 3.698816428990268
$ws.Range("I3").Value = 28.20978354406167
$ws.Range("J3").Value = 11.46585429578151

# Row 4 (index 2)
$ws.Range("B4").Value = 18.12066903969581
$ws.Range("C4").Value = 15.26374942226611
$ws.Range("D4").Value = 7.869580968758901
$ws.Range("F4").Value = 44.72686942835847
$ws.Range("G4").Value = 3.701723984364614
$ws.Range("I4").Value = 28.19752053530451
$ws.Range("J4").Value = 11.46241267123784

# Row 5 (index 3)
$ws.Range("B5").Value = 17.99690472198348
$ws.Range("C5").Value = 15.13311860654594
$ws.Range("D5").Value = 7.868608225489491
$ws.Range("F5").Value = 44.6624212837189
$ws.Range("G5").Value = 3.702943960910533
$ws.Range("I5").Value = 28.19450499335422
$ws.Range("J5").Value = 11.46165612716688

# Row 6 (index 4)
$ws.Range("B6").Value = 17.97635804810918
$ws.Range("C6").Value = 15.11139793282675
$ws.Range("D6").Value = 7.868465175521076
$ws.Range("F6").Value = 44.65192957450044
$ws.Range("G6").Value = 3.703148662917167
$ws.Range("I6").Value = 28.19412364720083
$ws.Range("J6").Value = 11.46156950156162

# Row 7 (index 5)
$ws.Range("B7").Value = 18.11899975923024
$ws.Range("C7").Value = 15.2619898225019
$ws.Range("D7").Value = 7.869566612920699
$ws.Range("F7").Value = 44.72598621106422
$ws.Range("G7").Value = 3.70174029497755
$ws.Range("I7").Value = 28.19747185635746
$ws.Range("J7").Value = 11.46239985343415

# Row 8 (index 6)
$ws.Range("B8").Value = 18.74761294161035
$ws.Range("C8").Value = 15.92060299290189
$ws.Range("D8").Value = 7.877097176778373
$ws.Range("F8").Value = 45.08303082971972
$ws.Range("G8").Value = 3.695833432395918
$ws.Range("I8").Value = 28.22967035357515
$ws.Range("J8").Value = 11.4717356580521

# Row 9 (index 7)
$ws.Range("B9").Value = 19.96417708372157
$ws.Range("C9").Value = 17.17667444458284
$ws.Range("D9").Value = 7.901310672292378
$ws.Range("F9").Value = 45.88896817776587
$ws.Range("G9").Value = 3.685346753915231
$ws.Range("I9").Value = 28.35598060409447
$ws.Range("J9").Value = 11.51042842587452

# Row 10 (index 8)
$ws.Range("B10").Value = 20.83735684119504
$ws.Range("C10").Value = 18.06670601399976
$ws.Range("D10").Value = 7.924715817780902
$ws.Range("F10").Value = 46.54237947773884
$ws.Range("G10").Value = 3.67830107950354
$ws.Range("I10").Value = 28.48764073838049
$ws.Range("J10").Value = 11.55125345987641

# Row 11 (index 9)
$ws.Range("B11").Value = 21.22806744095283
$ws.Range("C11").Value = 18.46245584547413
$ws.Range("D11").Value = 7.936549880488933
$ws.Range("F11").Value = 46.85214330352358
$ws.Range("G11").Value = 3.675236815815827
$ws.Range("I11").Value = 28.55600984073613
$ws.Range("J11").Value = 11.57250255412148

# Row 12 (index 10)
$ws.Range("B12").Value = 21.37493146533995
$ws.Range("C12").Value = 18.61085684022893
$ws.Range("D12").Value = 7.941198826040194
$ws.Range("F12").Value = 46.9711690590385
$ws.Range("G12").Value = 3.6740965518133
$ws.Range("I12").Value = 28.58311785722013
$ws.Range("J12").Value = 11.58093198142151

# Row 13 (index 11)
$ws.Range("B13").Value = 21.34335245090674
$ws.Range("C13").Value = 18.57896319135498
$ws.Range("D13").Value = 7.940190186044531
$ws.Range("F13").Value = 46.94545926984628
$ws.Range("G13").Value = 3.674341236193045
$ws.Range("I13").Value = 28.5772255172522
$ws.Range("J13").Value = 11.57909956508093

# Row 14 (index 12)
$ws.Range("B14").Value = 21.24017270818947
$ws.Range("C14").Value = 18.47469495393733
$ws.Range("D14").Value = 7.936929010027646
$ws.Range("F14").Value = 46.86190150394145
$ws.Range("G14").Value = 3.675142603499511
$ws.Range("I14").Value = 28.55821563187023
$ws.Range("J14").Value = 11.5731883891805

# Row 15 (index 13)
$ws.Range("B15").Value = 21.1768258770771
$ws.Range("C15").Value = 18.41063327823133
$ws.Range("D15").Value = 7.93495318550383
$ws.Range("F15").Value = 46.81094226934237
$ws.Range("G15").Value = 3.675636078142801
$ws.Range("I15").Value = 28.5467301077959
$ws.Range("J15").Value = 11.569617410432

# Row 16 (index 14)
$ws.Range("B16").Value = 20.81167862310864
$ws.Range("C16").Value = 18.04064609583353
$ws.Range("D16").Value = 7.923966061171478
$ws.Range("F16").Value = 46.52238129406132
$ws.Range("G16").Value = 3.678504157950837
$ws.Range("I16").Value = 28.48334320082168
$ws.Range("J16").Value = 11.54991848532456

# Row 17 (index 15)
$ws.Range("B17").Value = 20.58589294961096
$ws.Range("C17").Value = 17.81122211897397
$ws.Range("D17").Value = 7.917527695622977
$ws.Range("F17").Value = 46.34851341379754
$ws.Range("G17").Value = 3.680299601663347
$ws.Range("I17").Value = 28.44662844268648
$ws.Range("J17").Value = 11.53851821636934

# Row 18 (index 16)
$ws.Range("B18").Value = 20.45542329161684
$ws.Range("C18").Value = 17.67841196045125
$ws.Range("D18").Value = 7.913936428045868
$ws.Range("F18").Value = 46.24969283990302
$ws.Range("G18").Value = 3.681345559601539
$ws.Range("I18").Value = 28.4263090798669
$ws.Range("J18").Value = 11.53221326823739

# Row 19 (index 17)
$ws.Range("B19").Value = 20.41114973864503
$ws.Range("C19").Value = 17.63330324022532
$ws.Range("D19").Value = 7.912739797903128
$ws.Range("F19").Value = 46.21643942678003
$ws.Range("G19").Value = 3.681701985970987
$ws.Range("I19").Value = 28.41956633986357
$ws.Range("J19").Value = 11.53012189327647

# Row 20 (index 18)
$ws.Range("B20").Value = 20.60999186186366
$ws.Range("C20").Value = 17.83573394855056
$ws.Range("D20").Value = 7.918201510019697
$ws.Range("F20").Value = 46.36689998216607
$ws.Range("G20").Value = 3.680107101684185
$ws.Range("I20").Value = 28.45045420118293
$ws.Range("J20").Value = 11.53970570864908

# Row 21 (index 19)
$ws.Range("B21").Value = 21.27050983744641
$ws.Range("C21").Value = 18.50536181115476
$ws.Range("D21").Value = 7.937882371467005
$ws.Range("F21").Value = 46.88639825731783
$ws.Range("G21").Value = 3.674906677997874
$ws.Range("I21").Value = 28.56376625221538
$ws.Range("J21").Value = 11.57491427372447

# Row 22 (index 20)
$ws.Range("B22").Value = 21.69578300748098
$ws.Range("C22").Value = 18.93442962189338
$ws.Range("D22").Value = 7.951720615323405
$ws.Range("F22").Value = 47.23593655129201
$ws.Range("G22").Value = 3.67162502596976
$ws.Range("I22").Value = 28.64492064887276
$ws.Range("J22").Value = 11.6001553714678

# Row 23 (index 21)
$ws.Range("B23").Value = 21.46944121599453
$ws.Range("C23").Value = 18.7062572800284
$ws.Range("D23").Value = 7.944246643395148
$ws.Range("F23").Value = 47.04849103731014
$ws.Range("G23").Value = 3.673365837241015
$ws.Range("I23").Value = 28.60095826567749
$ws.Range("J23").Value = 11.58648048437739

# Row 24 (index 22)
$ws.Range("B24").Value = 20.59909879521031
$ws.Range("C24").Value = 17.824654977364
$ws.Range("D24").Value = 7.917896534946216
$ws.Range("F24").Value = 46.35858386307382
$ws.Range("G24").Value = 3.680194088136388
$ws.Range("I24").Value = 28.44872211991585
$ws.Range("J24").Value = 11.53916806696997

# Row 25 (index 23)
$ws.Range("B25").Value = 19.63795192920113
$ws.Range("C25").Value = 16.841891003694
$ws.Range("D25").Value = 7.893763294645977
$ws.Range("F25").Value = 45.88896817776587
$ws.Range("G25").Value = 3.688067272772163
$ws.Range("I25").Value = 28.31500391846392
$ws.Range("J25").Value = 11.51042842587452
